$d = $word.ActiveDocument
Write-Output "Hello"
